# Add carjacking data for 2022-08-13 (new day: August 05, 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-08-05"

# Update the header label in column B (the "through" date text)
$ws.Range("B1").Value = "August 2022 (through August 05)"

# New counts recorded for August 05, 2022 (column B = current "August" period)
# and incremental updates to other cells within the August-anchored columns
# (column B = Aug 2022, J = Aug 2021, R = Aug 2020, Z = Aug 2019,
#  AH = Aug 2018, AP = Aug 2017, AX = Aug 2016, BF = Aug 2015)

$ws.Range("AH4").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("AX7").Value = 1
$ws.Range("AX9").Value = 1
$ws.Range("BF9").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("B12").Value = 2
$ws.Range("AP12").Value = 2
$ws.Range("J13").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("BF15").Value = 2
$ws.Range("J16").Value = 1
$ws.Range("J18").Value = 1
$ws.Range("R18").Value = 1
$ws.Range("AP23").Value = 1
$ws.Range("B51").Value = 1
$ws.Range("R57").Value = 3
$ws.Range("R64").Value = 1
$ws.Range("J66").Value = 2
$ws.Range("Z69").Value = 1
$ws.Range("AP82").Value = 1
$ws.Range("R85").Value = 1
